$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.227.76"
$ws.Range("E2").Value = "  +3.89%  "
$ws.Range("D3").Value = "2.454.90"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.34%  "
$ws.Range("E7").Value = "  +0.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +4.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.47%  "
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").Value = "2.837.11"
$ws.Range("E15").Value = "  +1.46%  "
$ws.Range("D16").Value = "2.468.48"
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.843"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "46.067.54"
$ws.Range("E18").Value = "  +3.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.29%  "
$ws.Range("E23").Value = "  +4.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "248.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("E25").Value = "  +2.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "25.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.53%  "
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("E29").Value = "  +1.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "49.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("E32").Value = "  +5.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.68%  "
$ws.Range("E34").Value = "  +3.52%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0765"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.53%  "
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("E39").Value = "  +2.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "126.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.66%  "
$ws.Range("E41").Value = "  +1.87%  "
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("E44").Value = "  +1.47%  "
$ws.Range("D45").Value = "1.972.04"
$ws.Range("E45").Value = "  +1.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.21%  "
$ws.Range("E48").Value = "  +12.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.19%  "
$ws.Range("E50").Value = "  +9.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.02%  "
